$p = $ppt.ActivePresentation

# Remove the third slide (blank "Title 1" / "Content Placeholder 2" slide)
$s = $p.Slides.Item(3)
$s.Delete()
